# Auto-generated edit script: refresh market-price-derived columns (H-N)
# across multiple job sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2128.1304
$ws.Range("I137").Value = 2344.5557
$ws.Range("K137").Value = 7033.6671
$ws.Range("M137").Value = -4483.6671

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 6759696
$ws.Range("I138").Value = 1333.15
$ws.Range("J138").Value = 9262793
$ws.Range("K138").Value = 3999.45
$ws.Range("L138").Value = 27788379
$ws.Range("M138").Value = 1140.55
$ws.Range("N138").Value = -27798659

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 458.30768
$ws.Range("I4").Value = 282.8
$ws.Range("J4").Value = 1043.3334
$ws.Range("K4").Value = 282.8
$ws.Range("L4").Value = 1043.3334
$ws.Range("M4").Value = -166.8
$ws.Range("N4").Value = -1275.3334

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3952.4119
$ws.Range("I61").Value = 2857.1936
$ws.Range("K61").Value = 2857.1936
$ws.Range("M61").Value = -2645.1936

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 15496.214
$ws.Range("I74").Value = 2339.2222
$ws.Range("K74").Value = 2339.2222
$ws.Range("M74").Value = -1465.2222

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 15496.214
$ws.Range("I77").Value = 2339.2222
$ws.Range("K77").Value = 11696.111
$ws.Range("M77").Value = -7328.111000000001

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 3496.5
$ws.Range("I102").Value = 3122.6365
$ws.Range("J102").Value = 4319
$ws.Range("K102").Value = 3122.6365
$ws.Range("L102").Value = 4319
$ws.Range("M102").Value = -1500.6365
$ws.Range("N102").Value = -7563

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2809
$ws.Range("I132").Value = 2416.0557
$ws.Range("K132").Value = 7248.1671
$ws.Range("M132").Value = -4718.1671

# Row 135 (Leve Item ID 42016)
$ws.Range("H135").Value = 56000
$ws.Range("J135").Value = 56000
$ws.Range("L135").Value = 56000
$ws.Range("N135").Value = -66140

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3952.4119
$ws.Range("I136").Value = 2857.1936
$ws.Range("K136").Value = 8571.5808
$ws.Range("M136").Value = -6021.5808

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2226
$ws.Range("I86").Value = 1580.9375
$ws.Range("J86").Value = 5666.3335
$ws.Range("K86").Value = 1580.9375
$ws.Range("L86").Value = 5666.3335
$ws.Range("M86").Value = -457.9375
$ws.Range("N86").Value = -7912.3335

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2226
$ws.Range("I89").Value = 1580.9375
$ws.Range("J89").Value = 5666.3335
$ws.Range("K89").Value = 7904.6875
$ws.Range("L89").Value = 28331.6675
$ws.Range("M89").Value = -2288.6875
$ws.Range("N89").Value = -39563.6675

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 43 (Leve Item ID 18504)
$ws.Range("H43").Value = 60000
$ws.Range("J43").Value = 60000
$ws.Range("L43").Value = 60000
$ws.Range("N43").Value = -60368

# Row 53 (Leve Item ID 25632)
$ws.Range("H53").Value = 29999.5
$ws.Range("J53").Value = 29999.5
$ws.Range("L53").Value = 29999.5
$ws.Range("N53").Value = -31213.5

# Row 95 (Leve Item ID 18192)
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

# Row 96 (Leve Item ID 18193)
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

# Row 97 (Leve Item ID 19730)
$ws.Range("H97").Value = 35000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982

# Row 101 (Leve Item ID 18504)
$ws.Range("H101").Value = 60000
$ws.Range("J101").Value = 60000
$ws.Range("L101").Value = 60000
$ws.Range("N101").Value = -66490

# Row 104 (Leve Item ID 19749)
$ws.Range("H104").Value = 42500
$ws.Range("I104").Value = 40000
$ws.Range("K104").Value = 40000
$ws.Range("M104").Value = -37379

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1573.6
$ws.Range("I105").Value = 1629.75
$ws.Range("J105").Value = 1349
$ws.Range("K105").Value = 1629.75
$ws.Range("L105").Value = 1349
$ws.Range("M105").Value = 117.25
$ws.Range("N105").Value = -4843

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1285.9412
$ws.Range("I107").Value = 927.8461
$ws.Range("J107").Value = 2449.75
$ws.Range("K107").Value = 927.8461
$ws.Range("L107").Value = 2449.75
$ws.Range("M107").Value = 992.1539
$ws.Range("N107").Value = -6289.75

# Row 108 (Leve Item ID 27087)
$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2443.125
$ws.Range("I132").Value = 2260.7942
$ws.Range("K132").Value = 6782.382599999999
$ws.Range("M132").Value = -4252.382599999999

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 9268.432000000001
$ws.Range("I134").Value = 5263.3413
$ws.Range("K134").Value = 15790.0239
$ws.Range("M134").Value = -13255.0239

# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 600081.5
$ws.Range("J141").Value = 600081.5
$ws.Range("L141").Value = 600081.5
$ws.Range("N141").Value = -610441.5

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 526.86365
$ws.Range("J23").Value = 557.4375
$ws.Range("L23").Value = 1672.3125
$ws.Range("N23").Value = -2142.3125

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 934.4091
$ws.Range("I122").Value = 366.66666
$ws.Range("J122").Value = 1147.3125
$ws.Range("K122").Value = 3299.99994
$ws.Range("L122").Value = 10325.8125
$ws.Range("M122").Value = -849.9999399999997
$ws.Range("N122").Value = -15225.8125

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3499.6667
$ws.Range("I80").Value = 3500
$ws.Range("K80").Value = 3500
$ws.Range("M80").Value = -2502

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3499.6667
$ws.Range("I83").Value = 3500
$ws.Range("K83").Value = 17500
$ws.Range("M83").Value = -12508

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 22727834
$ws.Range("I102").Value = 555.025
$ws.Range("J102").Value = 250000620
$ws.Range("K102").Value = 555.025
$ws.Range("L102").Value = 250000620
$ws.Range("M102").Value = 1066.975
$ws.Range("N102").Value = -250003864

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 9474.290000000001
$ws.Range("J126").Value = 2832.6667
$ws.Range("L126").Value = 8498.000100000001
$ws.Range("N126").Value = -13438.0001

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3263.5
$ws.Range("I132").Value = 2793.5217
$ws.Range("K132").Value = 8380.5651
$ws.Range("M132").Value = -5850.5651

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3620.1785
$ws.Range("I40").Value = 3054.8262
$ws.Range("K40").Value = 3054.8262
$ws.Range("M40").Value = -2918.8262

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 5030.773
$ws.Range("I122").Value = 4467.6875
$ws.Range("K122").Value = 13403.0625
$ws.Range("M122").Value = -10953.0625

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 15 (Leve Item ID 2670)
$ws.Range("H15").Value = 14665.667
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 16998.5
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 16998.5
$ws.Range("M15").Value = -9712
$ws.Range("N15").Value = -17574.5

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2464.1765
$ws.Range("I132").Value = 2242.2144
$ws.Range("K132").Value = 6726.6432
$ws.Range("M132").Value = -4196.6432

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1772.1613
$ws.Range("I136").Value = 1561.4615
$ws.Range("K136").Value = 4684.3845
$ws.Range("M136").Value = -2134.3845
